$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update progress values: row 7 (alinea c) and row 8 (alinea 1a) are now fully done (100%)
$ws.Range("D7").Value = 100
$ws.Range("D8").Value = 100

# Row 8 is now assigned to Bernardo instead of the previous "-" placeholder formula
$ws.Range("E8").Value = "Bernardo"

# Update the active selection on the sheet
$ws.Range("I12").Select()
